$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Target cluster (D2) and numeric columns ---
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.520808
$ws.Range("H2").Value = 7.562424
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.558821666666667
$ws.Range("N2").Value = 7.676465
$ws.Range("O2").Value = 0.2156728774407755
$ws.Range("P2").Value = 0.2156728774407755
$ws.Range("Q2").Value = 6.450298127906668
$ws.Range("R2").Value = 58.05268315116
$ws.Range("S2").Value = 0.2156728774407755
$ws.Range("T2").Value = 0.2156728774407755

# --- Row 3: update Target cluster (D3) and numeric columns ---
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.520808
$ws.Range("H3").Value = 7.562424
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("N3").Value = 19.361954
$ws.Range("O3").Value = 0.5439806384912759
$ws.Range("P3").Value = 0.5439806384912759
$ws.Range("Q3").Value = 16.26925617961067
$ws.Range("R3").Value = 146.423305616496
$ws.Range("S3").Value = 0.5439806384912759
$ws.Range("T3").Value = 0.5439806384912759

# --- Row 4: brand-new row ---
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.520808
$ws.Range("H4").Value = 7.562424
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.851558333333334
$ws.Range("N4").Value = 8.554675000000001
$ws.Range("O4").Value = 0.2403464840679487
$ws.Range("P4").Value = 0.2403464840679487
$ws.Range("Q4").Value = 7.188231059133335
$ws.Range("R4").Value = 64.69407953220001
$ws.Range("S4").Value = 0.2403464840679487
$ws.Range("T4").Value = 0.2403464840679487
